$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "feather esp32 v2" pinout ---
# Header (row 2) - plain/default style, matches existing header row text cells
$ws.Range("I2").Value = "feather esp32 v2"

# Numeric pin values for rows 5-9 (VCC, GND, RST rows align with existing columns)
$ws.Range("I5").Value = 21
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 5
$ws.Range("I8").Value = 14
$ws.Range("I9").Value = 27

# Row 10 is textual ("rst") instead of numeric
$ws.Range("I10").Value = "rst"

$ws.Range("I11").Value = 13
$ws.Range("I12").Value = 22
$ws.Range("I13").Value = 20
$ws.Range("I14").Value = 15
$ws.Range("I15").Value = 12

# --- New rows 17-18: button / buzzer notes ---
$ws.Range("E17").Value = "button"
$ws.Range("G17").Value = 17
$ws.Range("I17").Value = 38

$ws.Range("E18").Value = "output?"
$ws.Range("G18").Value = 12
$ws.Range("D18").Value = "buzzer"

# Apply the "label" style (same font/fill as the bordered data cells, but
# without the border) used by the new E17/E18/G18 cells. Format E17 first,
# then propagate the resulting style to E18/G18 via a format-only paste so
# that all three cells share a single new style entry.
$labelSrc = $ws.Range("E17")
$labelSrc.Font.Name = "Microsoft YaHei"
$labelSrc.Font.Size = 10
$labelSrc.Font.Color = 3355443
$labelSrc.Interior.Color = 16777215
$labelSrc.HorizontalAlignment = -4108
$labelSrc.VerticalAlignment = -4108
$labelSrc.WrapText = $true
$labelSrc.Borders.LineStyle = -4142

$labelSrc.Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("E18").Value = "output?"
$ws.Range("G18").Value = 12

# --- View/window cosmetics ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("I6").Select()
